# "задачи 2.xlsx" update:
#  - remove several task rows (the list of tasks got pruned down to 7 items)
#  - add "Теги" (tags) text for task 1 and input/output examples for
#    tasks 1, 2 and 4 (columns C/D)
#  - the remaining rows shift up, the trailing blank rows are trimmed so the
#    sheet ends at row 41 instead of row 50

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: "Программа должна переводить число ..." -----------------------
$ws.Range("C2").Value = "Входные данные:`n700`nВыходные данные:`n0.7 "
$ws.Range("D2").Value = "С++, сложность F, 2 семестр, типы данных"

# --- Row 3: "Дано четырехзначное число ..." --------------------------------
$ws.Range("C3").Value = "Входные данные:`n5693`nВыходные данные:`n3965"

# --- Row 4: "Даны значения двух моментов времени ..." (unchanged content) -
# (nothing to change, it just shifts row position because rows above it
#  were removed - see the row-delete step further down)

# --- Row 5 (was row 11): "Определите наименьшее расстояние ..." -----------
$ws.Range("B5").Value = 'Определите наименьшее расстояние между двумя локальными максимумами последовательности натуральных чисел, завершающейся числом 0. Если в последовательности нет двух локальных максимумов, выведите число 0.'
$ws.Range("C5").Value = "Входные данные:`n3 4 5 1 6 2 1 3 5 0`nВыходные данные:`n1"

# --- Row 6 (was row 12): "Даны действительные коэффициенты a, b, c ..." ---
$ws.Range("B6").Value = 'Даны действительные коэффициенты a, b, c. Решите уравнение ax2 + bx + c = 0 и выведите все его корни.Если данное уравнение не имеет корней, выведите число 0. Если уравнение имеет один корень, выведите число 1, а затем этот корень. Если уравнение имеет два корня, выведите число 2, а затем два корня в порядке возрастания. Если уравнение имеет бесконечно много корней, выведите число 3.'

# --- Row 7 (was row 14): "Дан массив ..." ----------------------------------
$ws.Range("B7").Value = 'Дан массив. Выведите те его элементы, которые встречаются в массиве только один раз. Элементы нужно выводить в том порядке, в котором они встречаются в списке.'
$ws.Range("C7").Value = "8`n4 3 5 2 5 1 3 5 вывод 4 2 1"

# --- Row 8 (was row 15): "Циклически сдвиньте элементы списка вправо ..." -
$ws.Range("B8").Value = 'Циклически сдвиньте элементы списка вправо (A[0] переходит на место A[1], A[1] на место A[2], ..., последний элемент переходит на место A[0]).'

# --- Clear the now-superseded source cells (old rows 9-17 content) --------
$ws.Range("B9").Value = ""
$ws.Range("B10").Value = ""
$ws.Range("B11").Value = ""
$ws.Range("B12").Value = ""
$ws.Range("B13").Value = ""
$ws.Range("B14").Value = ""
$ws.Range("C14").Value = ""
$ws.Range("B15").Value = ""
$ws.Range("B16").Value = ""
$ws.Range("B17").Value = ""

# --- Row heights (tasks now sized to a uniform 100pt, except row 6) -------
for ($r = 2; $r -le 5; $r++) {
    $ws.Rows($r).RowHeight = 100
}
$ws.Rows(6).RowHeight = 113.25
for ($r = 7; $r -le 41; $r++) {
    $ws.Rows($r).RowHeight = 100
}

# --- Column widths ----------------------------------------------------------
$ws.Columns(4).ColumnWidth = 28.21484375

# --- Remove the trailing blank rows (old rows 41-49); the old bottom-border
#     row 50 shifts up to become the new last row, 41 ------------------------
[void]$ws.Rows("41:49").Delete()

# --- View / selection state --------------------------------------------------
[void]$ws.Rows(9).Select()
$excel.ActiveWindow.Zoom = 70
